$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 15.428109
$ws.Range("H2").Value = 46.284327
$ws.Range("I2").Value = 0.105145687357564
$ws.Range("J2").Value = 0.105145687357564
$ws.Range("M2").Value = 0.5373756666666667
$ws.Range("N2").Value = 1.612127
$ws.Range("O2").Value = 0.007472820128982582
$ws.Range("P2").Value = 0.007472820128982581
$ws.Range("Q2").Value = 8.290690359281001
$ws.Range("R2").Value = 74.61621323352901
$ws.Range("S2").Value = 0.0007857348089613136
$ws.Range("T2").Value = 0.0007857348089613134
$ws.Range("G3").Value = 15.428109
$ws.Range("H3").Value = 46.284327
$ws.Range("I3").Value = 0.105145687357564
$ws.Range("J3").Value = 0.105145687357564
$ws.Range("O3").Value = 0.1537223653287423
$ws.Range("P3").Value = 0.1537223653287423
$ws.Range("Q3").Value = 170.546662470035
$ws.Range("R3").Value = 1534.919962230315
$ws.Range("S3").Value = 0.01616324376472118
$ws.Range("T3").Value = 0.01616324376472117
$ws.Range("G4").Value = 15.428109
$ws.Range("H4").Value = 46.284327
$ws.Range("I4").Value = 0.105145687357564
$ws.Range("J4").Value = 0.105145687357564
$ws.Range("M4").Value = 30.561198
$ws.Range("N4").Value = 91.683594
$ws.Range("O4").Value = 0.4249882340167162
$ws.Range("P4").Value = 0.4249882340167161
$ws.Range("Q4").Value = 471.5014939145821
$ws.Range("R4").Value = 4243.513445231239
$ws.Range("S4").Value = 0.04468567998456488
$ws.Range("T4").Value = 0.04468567998456487
$ws.Range("G5").Value = 15.428109
$ws.Range("H5").Value = 46.284327
$ws.Range("I5").Value = 0.105145687357564
$ws.Range("J5").Value = 0.105145687357564
$ws.Range("M5").Value = 29.75783666666667
$ws.Range("N5").Value = 89.27351
$ws.Range("O5").Value = 0.4138165805255589
$ws.Range("P5").Value = 0.4138165805255589
$ws.Range("Q5").Value = 459.10714769753
$ws.Range("R5").Value = 4131.96432927777
$ws.Range("S5").Value = 0.04351102879931662
$ws.Range("T5").Value = 0.04351102879931661
$ws.Range("I6").Value = 0.4914986618531588
$ws.Range("J6").Value = 0.4914986618531588
$ws.Range("M6").Value = 0.5373756666666667
$ws.Range("N6").Value = 1.612127
$ws.Range("O6").Value = 0.007472820128982582
$ws.Range("P6").Value = 0.007472820128982581
$ws.Range("Q6").Value = 38.75444937240555
$ws.Range("R6").Value = 348.79004435165
$ws.Range("S6").Value = 0.003672881093664288
$ws.Range("T6").Value = 0.003672881093664288
$ws.Range("I7").Value = 0.4914986618531588
$ws.Range("J7").Value = 0.4914986618531588
$ws.Range("O7").Value = 0.1537223653287423
$ws.Range("P7").Value = 0.1537223653287423
$ws.Range("S7").Value = 0.07555433685597927
$ws.Range("T7").Value = 0.07555433685597926
$ws.Range("I8").Value = 0.4914986618531588
$ws.Range("J8").Value = 0.4914986618531588
$ws.Range("M8").Value = 30.561198
$ws.Range("N8").Value = 91.683594
$ws.Range("O8").Value = 0.4249882340167162
$ws.Range("P8").Value = 0.4249882340167161
$ws.Range("Q8").Value = 2204.0119680107
$ws.Range("R8").Value = 19836.1077120963
$ws.Range("S8").Value = 0.2088811483225531
$ws.Range("T8").Value = 0.2088811483225531
$ws.Range("I9").Value = 0.4914986618531588
$ws.Range("J9").Value = 0.4914986618531588
$ws.Range("M9").Value = 29.75783666666667
$ws.Range("N9").Value = 89.27351
$ws.Range("O9").Value = 0.4138165805255589
$ws.Range("P9").Value = 0.4138165805255589
$ws.Range("Q9").Value = 2146.075168762722
$ws.Range("R9").Value = 19314.6765188645
$ws.Range("S9").Value = 0.2033902955809621
$ws.Range("T9").Value = 0.2033902955809621
$ws.Range("G10").Value = 17.753286
$ws.Range("H10").Value = 53.25985799999999
$ws.Range("I10").Value = 0.120992239510715
$ws.Range("J10").Value = 0.120992239510715
$ws.Range("M10").Value = 0.5373756666666667
$ws.Range("N10").Value = 1.612127
$ws.Range("O10").Value = 0.007472820128982582
$ws.Range("P10").Value = 0.007472820128982581
$ws.Range("Q10").Value = 9.540183899774
$ws.Range("R10").Value = 85.86165509796599
$ws.Range("S10").Value = 0.0009041532428663526
$ws.Range("T10").Value = 0.0009041532428663523
$ws.Range("G11").Value = 17.753286
$ws.Range("H11").Value = 53.25985799999999
$ws.Range("I11").Value = 0.120992239510715
$ws.Range("J11").Value = 0.120992239510715
$ws.Range("O11").Value = 0.1537223653287423
$ws.Range("P11").Value = 0.1537223653287423
$ws.Range("Q11").Value = 196.24982395289
$ws.Range("R11").Value = 1766.24841557601
$ws.Range("S11").Value = 0.01859921324400882
$ws.Range("T11").Value = 0.01859921324400881
$ws.Range("G12").Value = 17.753286
$ws.Range("H12").Value = 53.25985799999999
$ws.Range("I12").Value = 0.120992239510715
$ws.Range("J12").Value = 0.120992239510715
$ws.Range("M12").Value = 30.561198
$ws.Range("N12").Value = 91.683594
$ws.Range("O12").Value = 0.4249882340167162
$ws.Range("P12").Value = 0.4249882340167161
$ws.Range("Q12").Value = 542.561688596628
$ws.Range("R12").Value = 4883.055197369651
$ws.Range("S12").Value = 0.05142027819938631
$ws.Range("T12").Value = 0.0514202781993863
$ws.Range("G13").Value = 17.753286
$ws.Range("H13").Value = 53.25985799999999
$ws.Range("I13").Value = 0.120992239510715
$ws.Range("J13").Value = 0.120992239510715
$ws.Range("M13").Value = 29.75783666666667
$ws.Range("N13").Value = 89.27351
$ws.Range("O13").Value = 0.4138165805255589
$ws.Range("P13").Value = 0.4138165805255589
$ws.Range("Q13").Value = 528.29938508462
$ws.Range("R13").Value = 4754.69446576158
$ws.Range("S13").Value = 0.05006859482445349
$ws.Range("T13").Value = 0.05006859482445349
$ws.Range("G14").Value = 41.43140433333334
$ws.Range("H14").Value = 124.294213
$ws.Range("I14").Value = 0.2823634112785623
$ws.Range("J14").Value = 0.2823634112785622
$ws.Range("M14").Value = 0.5373756666666667
$ws.Range("N14").Value = 1.612127
$ws.Range("O14").Value = 0.007472820128982582
$ws.Range("P14").Value = 0.007472820128982581
$ws.Range("Q14").Value = 22.26422852456123
$ws.Range("R14").Value = 200.378056721051
$ws.Range("S14").Value = 0.002110050983490628
$ws.Range("T14").Value = 0.002110050983490627
$ws.Range("G15").Value = 41.43140433333334
$ws.Range("H15").Value = 124.294213
$ws.Range("I15").Value = 0.2823634112785623
$ws.Range("J15").Value = 0.2823634112785622
$ws.Range("O15").Value = 0.1537223653287423
$ws.Range("P15").Value = 0.1537223653287423
$ws.Range("Q15").Value = 457.9944133462207
$ws.Range("R15").Value = 4121.949720115986
$ws.Range("S15").Value = 0.04340557146403308
$ws.Range("T15").Value = 0.04340557146403306
$ws.Range("G16").Value = 41.43140433333334
$ws.Range("H16").Value = 124.294213
$ws.Range("I16").Value = 0.2823634112785623
$ws.Range("J16").Value = 0.2823634112785622
$ws.Range("M16").Value = 30.561198
$ws.Range("N16").Value = 91.683594
$ws.Range("O16").Value = 0.4249882340167162
$ws.Range("P16").Value = 0.4249882340167161
$ws.Range("Q16").Value = 1266.193351249058
$ws.Range("R16").Value = 11395.74016124152
$ws.Range("S16").Value = 0.1200011275102119
$ws.Range("T16").Value = 0.1200011275102119
$ws.Range("G17").Value = 41.43140433333334
$ws.Range("H17").Value = 124.294213
$ws.Range("I17").Value = 0.2823634112785623
$ws.Range("J17").Value = 0.2823634112785622
$ws.Range("M17").Value = 29.75783666666667
$ws.Range("N17").Value = 89.27351
$ws.Range("O17").Value = 0.4138165805255589
$ws.Range("P17").Value = 0.4138165805255589
$ws.Range("Q17").Value = 1232.908963021959
$ws.Range("R17").Value = 11096.18066719763
$ws.Range("S17").Value = 0.1168466613208267
$ws.Range("T17").Value = 0.1168466613208267
